$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03449230194091797
$ws.Range("C2").Value = 0.06408958435058594
$ws.Range("D2").Value = 0.008924674987792969
$ws.Range("E2").Value = 0.04052743911743164
$ws.Range("F2").Value = 0.00169215202331543
$ws.Range("G2").Value = 0.1265068054199219
$ws.Range("H2").Value = 0.0276667594909668
$ws.Range("I2").Value = 0.04682331085205078
$ws.Range("J2").Value = 0.01683006286621094
$ws.Range("K2").Value = 0.03265647888183594
$ws.Range("L2").Value = 0.008340692520141602
$ws.Range("M2").Value = 0.05038094520568848
$ws.Range("B3").Value = 0.1777416706085205
$ws.Range("C3").Value = 0.06609849929809571
$ws.Range("D3").Value = 0.04581508636474609
$ws.Range("E3").Value = 0.03242983818054199
$ws.Range("F3").Value = 0.02811088562011719
$ws.Range("G3").Value = 0.02763156890869141
$ws.Range("H3").Value = 0.2230873107910156
$ws.Range("I3").Value = 0.0617856502532959
$ws.Range("J3").Value = 0.1572887897491455
$ws.Range("K3").Value = 0.04955000877380371
$ws.Range("L3").Value = 0.05497927665710449
$ws.Range("M3").Value = 0.02655735015869141
$ws.Range("B4").Value = 0.06199078559875489
$ws.Range("C4").Value = 0.03156571388244629
$ws.Range("D4").Value = 0.02436656951904297
$ws.Range("E4").Value = 0.01900358200073242
$ws.Range("F4").Value = 0.1371804714202881
$ws.Range("G4").Value = 0.01466960906982422
$ws.Range("H4").Value = 0.04345941543579102
$ws.Range("I4").Value = 0.02898030281066894
$ws.Range("J4").Value = 0.03531708717346192
$ws.Range("K4").Value = 0.02787270545959473
$ws.Range("L4").Value = 0.08633460998535156
$ws.Range("M4").Value = 0.0222536563873291
$ws.Range("B5").Value = 0.03932452201843262
$ws.Range("C5").Value = 0.03577046394348145
$ws.Range("D5").Value = 0.03139300346374511
$ws.Range("E5").Value = 0.02863397598266602
$ws.Range("H5").Value = 0.0320746898651123
$ws.Range("I5").Value = 0.03615107536315918
$ws.Range("J5").Value = 0.02557921409606934
$ws.Range("K5").Value = 0.03418269157409668
$ws.Range("B6").Value = 0.4865874767303467
$ws.Range("C6").Value = 0.1022066593170166
$ws.Range("D6").Value = 0.4001742839813233
$ws.Range("E6").Value = 0.09559688568115235
$ws.Range("F6").Value = 0.09910163879394532
$ws.Range("G6").Value = 0.02847757339477539
$ws.Range("H6").Value = 0.5101908206939697
$ws.Range("I6").Value = 0.08968691825866699
$ws.Range("J6").Value = 0.3879057884216309
$ws.Range("K6").Value = 0.06888542175292969
$ws.Range("L6").Value = 0.1584257125854492
$ws.Range("M6").Value = 0.0247398853302002
